$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 46; this pushes the existing rows 46-64 down to 47-65
$ws.Rows.Item(46).Insert()

# Populate the new row 46 with the new record.
# Static columns match the surrounding rows (A, B, C, E, F, G, H, I, N, O, Q, R)
$ws.Cells.Item(46, 1).Value = 9
$ws.Cells.Item(46, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(46, 3).Value = "Metropolitana"
$ws.Cells.Item(46, 4).Value = 44806
$ws.Cells.Item(46, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(46, 5).Value = 13
$ws.Cells.Item(46, 6).Value = 100112035
$ws.Cells.Item(46, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(46, 8).Value = "Sin especificar"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 27
$ws.Cells.Item(46, 11).Value = 19000
$ws.Cells.Item(46, 12).Value = 20000
$ws.Cells.Item(46, 13).Value = 19556
$ws.Cells.Item(46, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(46, 15).Value = "Hijuelas"
$ws.Cells.Item(46, 16).Value = 1304
$ws.Cells.Item(46, 17).Value = 15
$ws.Cells.Item(46, 18).Value = "Hortaliza"
